# #12 Explanation for pyinstaller added on solution slides
#
# Slide 13 ("Pyinstaller") has a content placeholder whose first paragraph
# reads "We choose to use pyinstaller to deploy SAI on Windows and OS X",
# followed by several empty paragraphs. We insert a brand new paragraph
# right after that first paragraph (and before the empty ones) explaining
# that Pyinstaller lets us build an executable for every platform.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The paragraph that currently ends with "...deploy SAI on Windows and OS X".
$firstPara = $tr.Paragraphs(1, 1)

$cr = [char]13

# Each element becomes its own run, matching how PowerPoint splits runs
# along word boundaries as the text is typed/autocorrected.
$runs = @(
    "Pyinstaller",
    " ",
    "allows",
    " us to ",
    "create",
    " ",
    "executable",
    " for ",
    "all platforms"
)

# Create the new paragraph (prefixed with a paragraph break) holding the
# first run, inserted immediately after the first paragraph.
$null = $firstPara.InsertAfter($cr + $runs[0])

# Append the remaining runs to that freshly created second paragraph, one
# at a time so each becomes its own <a:r> run.
for ($i = 1; $i -lt $runs.Length; $i++) {
    $full = $shp.TextFrame.TextRange
    $newPara = $full.Paragraphs(2, 1)
    $null = $newPara.InsertAfter($runs[$i])
}
